$wb = $excel.ActiveWorkbook

# --- Sheet 1: LoginWithInvalidCredentialsTest ---
$ws1 = $wb.Worksheets.Item(1)

# Move the email/password values from row 5 up to row 2
$ws1.Range("D2").Value = "wiasm.mtour@gmail.com"
$ws1.Range("E2").Value = 12345678
$ws1.Range("D5").ClearContents()
$ws1.Range("E5").ClearContents()

# Update selection (was C3) -> D5:E5 with active cell D5
$ws1.Range("D5:E5").Select()

# --- Sheet 2: LoginWithValidCredentialsTest ---
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("E2").Value = 12345678955
$ws2.Range("E3").Value = 123456789
$ws2.Range("E4").Value = 123456789
$ws2.Range("E5").Value = 123456789

# Update selection (was E5) -> E2, and keep sheet2 as the active tab
$ws2.Range("E2").Select()
